$wb = $excel.ActiveWorkbook

# --- "Prepay Loan" sheet: add two "wait to page load" steps ---
$ws = $wb.Worksheets.Item("Prepay Loan")

# 1) Insert "waittopageload" / 4000 as the new row 2 (right after "clickonPrepayLoan").
#    Excel inherits the new blank row's formatting from the row above it, which already
#    matches the style we need (s="15" on both columns).
$ws.Rows(2).Insert()
$ws.Cells.Item(2, 1).Value = "waittopageload"
$ws.Cells.Item(2, 2).Value = 4000

# 2) Insert "waittopageload2" / 4000 right before the final "clickonsubmit" row.
#    To get the correct cell style (matching the "clickonsubmit" row, not the amount row
#    above it), insert the blank row AFTER "clickonsubmit" instead (so it inherits that
#    row's formatting), then shift the values: "clickonsubmit" moves down one row and the
#    new "waittopageload2" row takes its old slot.
$ws.Rows(7).Insert()
$ws.Cells.Item(7, 1).Value = $ws.Cells.Item(6, 1).Value2
$ws.Cells.Item(7, 2).Value = $ws.Cells.Item(6, 2).Value2
$ws.Cells.Item(6, 1).Value = "waittopageload2"
$ws.Cells.Item(6, 2).Value = 4000

# --- Update the remembered cell selection on each sheet ---
$ws.Activate()
$ws.Range("D13").Select()

$wb.Worksheets.Item("Summary").Activate()
$wb.Worksheets.Item("Summary").Range("A4").Select()

$wb.Worksheets.Item("Repayment schedule").Activate()
$wb.Worksheets.Item("Repayment schedule").Range("I3").Select()

$wb.Worksheets.Item("Transactions").Activate()
$wb.Worksheets.Item("Transactions").Range("H2").Select()
